# Apply the scene-cat block-order swap edit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the header labels for columns B and E
$ws.Range("B1").Value = "living_rooms_2"
$ws.Range("E1").Value = "kitchens_1"

# Update the one-hot indicator values so that each column keeps
# tracking the same semantic category after the header swap.
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 1

$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0

$ws.Range("C5").Value = 1
$ws.Range("D5").Value = 0

$ws.Range("B6").Value = 0
$ws.Range("F6").Value = 1

$ws.Range("B7").Value = 1
$ws.Range("E7").Value = 0
